$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) mirroring the existing header style used by the
# other header cells in row 1 (bold, centered, bordered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data rows for the new column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
